$d = $word.ActiveDocument
$d.Content.Find.Execute("Software engineer", $true, $false, $false, $false, $false, $true, 1, $false, "Senior System Architect", 2)
